$wb = $excel.ActiveWorkbook

# --- Update the "Logs" sheet: append a new row of test-mail data ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A11").Value = "Staan er nog EcoPro-700 op voorraad?"
$logs.Range("B11").Value = "mailmind.test@zohomail.eu"
$logs.Range("C11").Value = "Testmail #1: Staan er nog EcoPro-700 op voorraad?"
$logs.Range("D11").Value = "Inkoop / Bestellingen"
$logs.Range("E11").Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Range("F11").Value = "2025-08-06 19:54:49"
$logs.Range("G11").Value = "Ja"
$logs.Range("H11").Value = "Ja"
$logs.Range("I11").Value = "Nee"
$logs.Range("J11").Value = "Nee"

# --- Update the "Dashboard" sheet: bump the "Inkoop / Bestellingen" count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 4

# --- Extend the conditional formatting ranges to cover the new row ---
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "10")
    $newRange = $logs.Range($col + "2:" + $col + "11")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
